$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the CasesTab query text in B2: drop the erroneous trailing
#     `coalesce(co.cohort_description, '') AS `Cohort`` line (and its
#     trailing newline) that was producing a query error. ---
$newQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" + `
            "WHERE demo.breed IN ['Staffordshire Bull Terrier'] `n" + `
            "MATCH (c)<--(diag:diagnosis)`n" + `
            "OPTIONAL MATCH (samp:sample)-->(c)`n" + `
            "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" + `
            "WITH DISTINCT c, s, demo, diag, co`n" + `
            "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" + `
            "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" + `
            "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" + `
            "        coalesce(demo.breed, '') AS Breed ,`n" + `
            "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" + `
            "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" + `
            "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" + `
            "        coalesce(demo.sex, '') AS Sex ,`n" + `
            "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" + `
            "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" + `
            "        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value2 = $newQuery

# --- Row heights shrink a bit now that the text in B2 (and the
#     re-rendered B3/B4) take up slightly less vertical space. ---
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 216
$ws.Rows.Item(4).RowHeight = 244.8

# --- View state: selection moves to B2, zoom changes from 40% to 70%,
#     and the window no longer forces topLeftCell to A4. ---
[void]$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 70
